$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column O data for year 2021
$ws.Range("O4").Value = 2021
$ws.Range("O4").Style = $ws.Range("N4").Style

$ws.Range("O5").Value = 515
$ws.Range("O5").Style = $ws.Range("N5").Style

# Update the view: clear frozen/scrolled topLeftCell and update selection
$ws.Range("P12").Select()
